$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row 3: IPA111 test case (set in this order so new shared strings are
# appended in the same order as the target workbook)
$ws.Range("A3").Value = "IPA111"
$ws.Range("B3").Value = "OBT"
$ws.Range("C3").Value = "Save the company search data and rerun the saved data"

# Update row 2, column C: new description for the existing IPA001 test case
$ws.Range("C2").Value = "Save the technology search data and rerun the saved data"

# Apply the thin border (matching existing border used elsewhere) without any fill
$newRow = $ws.Range("A3:E3")
$newRow.Borders.Color = 0
$newRow.Borders.LineStyle = 1

# Update the active selection to reflect the newly added row
$null = $ws.Range("A3:E3").Select()
